$wb = $excel.ActiveWorkbook

# 展览 (28 changes)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 27   # was 26
$ws.Range("F3").Value = 105   # was 104
$ws.Range("F8").Value = 62   # was 59
$ws.Range("F9").Value = 6   # was 5
$ws.Range("F10").Value = 10012   # was 9975
$ws.Range("F14").Value = 121   # was 120
$ws.Range("F15").Value = 1937   # was 1930
$ws.Range("F16").Value = 871   # was 867
$ws.Range("F20").Value = 149   # was 147
$ws.Range("F22").Value = 215   # was 213
$ws.Range("F23").Value = 1090   # was 1086
$ws.Range("F24").Value = 65   # was 62
$ws.Range("F25").Value = 100   # was 99
$ws.Range("F28").Value = 127   # was 123
$ws.Range("F29").Value = 604   # was 600
$ws.Range("F30").Value = 2664   # was 2650
$ws.Range("F31").Value = 922   # was 918
$ws.Range("F32").Value = 606   # was 603
$ws.Range("F36").Value = 468   # was 458
$ws.Range("F39").Value = 1180   # was 1177
$ws.Range("F40").Value = 192   # was 184
$ws.Range("F41").Value = 101   # was 102
$ws.Range("F42").Value = 54   # was 53
$ws.Range("F43").Value = 101   # was 100
$ws.Range("F44").Value = 104   # was 94
$ws.Range("F45").Value = 24   # was 23
$ws.Range("F46").Value = 4033   # was 4032
$ws.Range("F47").Value = 44   # was 43

# 演出 (4 changes)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 1573   # was 27
$ws.Range("F7").Value = 38   # was 30
$ws.Range("F9").Value = 31   # was 30
$ws.Range("F11").Value = 100   # was 99

# 本地生活 (1 changes)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2355   # was 2354

# 全部类型 (27 changes)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 27   # was 26
$ws.Range("F5").Value = 105   # was 104
$ws.Range("F9").Value = 62   # was 59
$ws.Range("F10").Value = 6   # was 5
$ws.Range("F11").Value = 10012   # was 9976
$ws.Range("F15").Value = 121   # was 120
$ws.Range("F16").Value = 1937   # was 1930
$ws.Range("F17").Value = 871   # was 867
$ws.Range("F20").Value = 149   # was 147
$ws.Range("F22").Value = 215   # was 213
$ws.Range("F23").Value = 1090   # was 1086
$ws.Range("F24").Value = 65   # was 62
$ws.Range("F25").Value = 100   # was 99
$ws.Range("F26").Value = 1573   # was 27
$ws.Range("F30").Value = 127   # was 123
$ws.Range("F31").Value = 604   # was 600
$ws.Range("F32").Value = 2664   # was 2650
$ws.Range("F33").Value = 922   # was 918
$ws.Range("F34").Value = 38   # was 30
$ws.Range("F36").Value = 606   # was 603
$ws.Range("F38").Value = 468   # was 459
$ws.Range("F40").Value = 54   # was 53
$ws.Range("F41").Value = 101   # was 100
$ws.Range("F42").Value = 105   # was 94
$ws.Range("F43").Value = 24   # was 23
$ws.Range("F44").Value = 4033   # was 4032
$ws.Range("F48").Value = 44   # was 43
